# Recipient Heat Fuel Fractions.xlsx -- integrate EU data for RHFF
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet -- rewrite the narrative notes to describe the EU assumptions
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A1").Value = "RHFF Recipient Heat Fuel Fractions"
$about.Range("A3").Value = "Source:"
$about.Range("B3").Value = "none"
$about.Range("A5").Value = "Notes:"
$about.Range("A6").Value = "This variable specifies the recipient fuel for the Fraction of District Heat Fuel Use Shifted"
$about.Range("A7").Value = "to Other Fuels policy."
$about.Range("A8").Value = ""
$about.Range("A9").Value = "We assume that district heating in the EU will be decarbonized "
$about.Range("A10").Value = "through a shift from fossils to a mix of large scale heat pumps, "
$about.Range("A11").Value = "biomass and hydrogen. We further expect that oil and coal "
$about.Range("A12").Value = "retirements in the near future will, to a certain degree, be replaced "
$about.Range("A13").Value = "fossil fuels with a lower carbon intensity, such as natural gas, as this is"
$about.Range("A14").Value = "the picture we are currently seeing in the EU. "

# match the recorded post-edit selection on the About sheet
$about.Range("A15").Select()

# ---------------------------------------------------------------------------
# "RHFF" sheet -- update the recipient-fuel fraction matrix with EU figures
# ---------------------------------------------------------------------------
$rhff = $wb.Worksheets.Item("RHFF")

# Row 2: electricity (from-type)
$rhff.Range("B2").Value = 1
$rhff.Range("C2").Value = 0.25
$rhff.Range("D2").Value = 0.4
$rhff.Range("E2").Value = 0.5
$rhff.Range("F2").Value = 0.25
$rhff.Range("G2").Value = 0.5
$rhff.Range("H2").Value = 0.25
$rhff.Range("I2").Value = 0.25
$rhff.Range("J2").Value = 0.25
$rhff.Range("K2").Value = 0.4

# Row 4: natural gas (from-type)
$rhff.Range("C4").Value = 0.25
$rhff.Range("F4").Value = 0.25
$rhff.Range("H4").Value = 0.25
$rhff.Range("I4").Value = 0.25
$rhff.Range("J4").Value = 0.25

# Row 5: biomass (from-type)
$rhff.Range("C5").Value = 0.25
$rhff.Range("D5").Value = 0.2
$rhff.Range("F5").Value = 0.25
$rhff.Range("H5").Value = 0.25
$rhff.Range("I5").Value = 0.25
$rhff.Range("J5").Value = 0.25
$rhff.Range("K5").Value = 0.2

# Row 11: hydrogen (from-type) -- was all 1s, now mirrors the electricity row
$rhff.Range("B11").Value = 0
$rhff.Range("C11").Value = 0.25
$rhff.Range("D11").Value = 0.4
$rhff.Range("E11").Value = 0.5
$rhff.Range("F11").Value = 0.25
$rhff.Range("G11").Value = 0.5
$rhff.Range("H11").Value = 0.25
$rhff.Range("I11").Value = 0.25
$rhff.Range("J11").Value = 0.25
$rhff.Range("K11").Value = 0.4

# match the recorded post-edit selection on the RHFF sheet
$rhff.Range("B11").Select()
